# Reorder the player table so that "Herbert Jones" (previously data row 8)
# and "Malik Monk" (previously data row 12) become the first two data rows,
# with every other player keeping their existing relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two rows that are moving to the top of the list
# (right under the header row).
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# The inserted rows can pick up formatting from the row above (the bold
# header) — reset them back to the plain/default look used by every other
# data row before filling in their values.
$ws.Range("A2:C3").ClearFormats()

# Row 2: Herbert Jones
$ws.Cells.Item(2, 1).Value = "Herbert Jones"
$ws.Cells.Item(2, 2).Value = "SF,PF"
$ws.Cells.Item(2, 3).Value = "New Orleans Pelicans"

# Row 3: Malik Monk
$ws.Cells.Item(3, 1).Value = "Malik Monk"
$ws.Cells.Item(3, 2).Value = "PG,SG,SF"
$ws.Cells.Item(3, 3).Value = "Sacramento Kings"

# Remove the two now-duplicated rows further down the sheet (their former
# positions have shifted down by 2 because of the inserts above: row 8 ->
# row 10 for Herbert Jones, row 12 -> row 14 for Malik Monk).
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(13).Delete()
